# Update RAM labels on rows 16 and 18 (drop the "L" in "DDR3L" to match "DDR3")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "DDR3 2GB"
$ws.Range("E18").Value = "DDR3 1GB"

# Update the view: zoom level and active selection
$excel.ActiveWindow.Zoom = 175
$ws.Range("M32").Select()
